$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make row 6 match the styling already used by the rows above it (date /
# hours formats) before filling in the values, so no new styles get
# fabricated by assigning a raw .NET DateTime.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in row 6 with a new progress report entry
$ws.Range("A6").Value = "2/5/2017"
$ws.Range("B6").Formula = "=0.3"
$ws.Range("C6").Value = "Revising Phase I Document"

# Move the active selection to A7, as in the source file
$ws.Range("A7").Select()
